$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ON" / 5.96 / "[O-]N" data row (originally row 23) entirely.
# This shifts all subsequent rows up by one.
$ws.Rows.Item(23).Delete()

# Update the active selection to match the post-edit state.
$ws.Range("D27").Select()
